$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray header scratch values in H1:P1, shifting remaining
# cells left so the used range collapses back down to columns A:D.
$ws.Range("H1:P1").Delete(-4159)

# Prime the date format for the new rows by copying the existing
# date-formatted cell's format (reuses the same style as A2:A49
# instead of creating a new numFmtId).
$ws.Range("A49").Copy()
$ws.Range("A50:A73").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New temperature log entries collected on 2017-08-02 (serial 42949).
# Column C uses "N/A" because no tote reading was recorded this round.
$newData = @(
  @(1, 5),
  @(2, 6),
  @(3, 4),
  @(4, 3),
  @(5, 9),
  @(6, 7),
  @(7, 2),
  @(8, 4),
  @(9, 3),
  @(10, 8),
  @(11, 2),
  @(12, 6),
  @(13, 6),
  @(14, 2),
  @(15, 9),
  @(16, 5),
  @(17, 7),
  @(18, 8),
  @(19, 7),
  @(20, 8),
  @(21, 4),
  @(22, 3),
  @(23, 1),
  @(24, 9)
)

$row = 50
foreach ($item in $newData) {
    $bucket = $item[0]
    $tote = $item[1]
    $ws.Cells.Item($row, 1).Value = 42949
    $ws.Cells.Item($row, 2).Value = $bucket
    $ws.Cells.Item($row, 3).Value = "N/A"
    $ws.Cells.Item($row, 4).Value = $tote
    $row = $row + 1
}

# Match the author's final selection position.
$null = $ws.Range("J19").Select()
